# Auto-generated script to apply numeric cell updates per the commit diff.
# Each (sheet, row, column) triplet below sets the cell to its new value;
# a value of $null / empty string clears (removes) the cell, matching cells
# that were deleted in the diff (e.g. CUL N68, N71); newly introduced cells
# (e.g. CUL M81, M84, WVR N113) are simply set for the first time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 992.7742
$ws.Range("I74").Value = 947.087
$ws.Range("J74").Value = 1124.125
$ws.Range("K74").Value = 947.087
$ws.Range("L74").Value = 1124.125
$ws.Range("M74").Value = -73.08699999999999
$ws.Range("N74").Value = -2872.125
$ws.Range("H77").Value = 992.7742
$ws.Range("I77").Value = 947.087
$ws.Range("J77").Value = 1124.125
$ws.Range("K77").Value = 4735.434999999999
$ws.Range("L77").Value = 5620.625
$ws.Range("M77").Value = -367.4349999999995
$ws.Range("N77").Value = -14356.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1968.8125
$ws.Range("I31").Value = 1149.725
$ws.Range("J31").Value = 6064.25
$ws.Range("K31").Value = 1149.725
$ws.Range("L31").Value = 6064.25
$ws.Range("M31").Value = -854.7249999999999
$ws.Range("N31").Value = -6654.25
$ws.Range("H34").Value = 1968.8125
$ws.Range("I34").Value = 1149.725
$ws.Range("J34").Value = 6064.25
$ws.Range("K34").Value = 1149.725
$ws.Range("L34").Value = 6064.25
$ws.Range("M34").Value = -947.7249999999999
$ws.Range("N34").Value = -6468.25
$ws.Range("H132").Value = 1626.5312
$ws.Range("I132").Value = 1185
$ws.Range("J132").Value = 2271.8462
$ws.Range("K132").Value = 3555
$ws.Range("L132").Value = 6815.5386
$ws.Range("M132").Value = -1025
$ws.Range("N132").Value = -11875.5386
$ws.Range("H134").Value = 1379.6666
$ws.Range("I134").Value = 1200.2188
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 3600.6564
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -1065.6564
$ws.Range("N134").Value = -11670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 907.2162
$ws.Range("J5").Value = 3475
$ws.Range("L5").Value = 10425
$ws.Range("N5").Value = -10649
$ws.Range("H68").Value = 25575.5
$ws.Range("I68").Value = 25575.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 76726.5
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -75915.5
$ws.Range("N68").Value = ""
$ws.Range("H69").Value = 918.9091
$ws.Range("I69").Value = 702.6667
$ws.Range("J69").Value = 1000
$ws.Range("K69").Value = 2108.0001
$ws.Range("L69").Value = 3000
$ws.Range("M69").Value = -1297.0001
$ws.Range("N69").Value = -4622
$ws.Range("H71").Value = 25575.5
$ws.Range("I71").Value = 25575.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 230179.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -226123.5
$ws.Range("N71").Value = ""
$ws.Range("H72").Value = 918.9091
$ws.Range("I72").Value = 702.6667
$ws.Range("J72").Value = 1000
$ws.Range("K72").Value = 6324.0003
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = -2268.0003
$ws.Range("N72").Value = -17112
$ws.Range("H80").Value = 3276.8462
$ws.Range("I80").Value = 1399
$ws.Range("J80").Value = 3433.3333
$ws.Range("K80").Value = 4197
$ws.Range("L80").Value = 10299.9999
$ws.Range("M80").Value = -3261
$ws.Range("N80").Value = -12171.9999
$ws.Range("H81").Value = 1150
$ws.Range("I81").Value = 300
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 900
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = 223
$ws.Range("N81").Value = -8246
$ws.Range("H83").Value = 3276.8462
$ws.Range("I83").Value = 1399
$ws.Range("J83").Value = 3433.3333
$ws.Range("K83").Value = 12591
$ws.Range("L83").Value = 30899.9997
$ws.Range("M83").Value = -7911
$ws.Range("N83").Value = -40259.9997
$ws.Range("H84").Value = 1150
$ws.Range("I84").Value = 300
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 2700
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = 2916
$ws.Range("N84").Value = -29232
$ws.Range("H113").Value = 6211587
$ws.Range("I113").Value = 414.76923
$ws.Range("J113").Value = 14286111
$ws.Range("K113").Value = 1244.30769
$ws.Range("L113").Value = 42858333
$ws.Range("M113").Value = 925.6923099999999
$ws.Range("N113").Value = -42862673
$ws.Range("H131").Value = 938.6957
$ws.Range("I131").Value = 564.3333
$ws.Range("J131").Value = 994.85
$ws.Range("K131").Value = 1692.9999
$ws.Range("L131").Value = 2984.55
$ws.Range("M131").Value = 3347.0001
$ws.Range("N131").Value = -13064.55
$ws.Range("H135").Value = 907.2162
$ws.Range("J135").Value = 3475
$ws.Range("L135").Value = 31275
$ws.Range("N135").Value = -36345

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9433.923000000001
$ws.Range("I113").Value = 15663
$ws.Range("J113").Value = 2166.6667
$ws.Range("K113").Value = 15663
$ws.Range("L113").Value = 2166.6667
$ws.Range("M113").Value = -13493
$ws.Range("N113").Value = -6506.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 380.57895
$ws.Range("I22").Value = 317.5625
$ws.Range("J22").Value = 716.6667
$ws.Range("K22").Value = 317.5625
$ws.Range("L22").Value = 716.6667
$ws.Range("M22").Value = -22.5625
$ws.Range("N22").Value = -1306.6667
$ws.Range("H27").Value = 380.57895
$ws.Range("I27").Value = 317.5625
$ws.Range("J27").Value = 716.6667
$ws.Range("K27").Value = 317.5625
$ws.Range("L27").Value = 716.6667
$ws.Range("M27").Value = -210.5625
$ws.Range("N27").Value = -930.6667
$ws.Range("H61").Value = 4126
$ws.Range("I61").Value = 3752
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 3752
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -3550
$ws.Range("N61").Value = -4904
$ws.Range("H97").Value = 21367.5
$ws.Range("J97").Value = 21367.5
$ws.Range("L97").Value = 21367.5
$ws.Range("N97").Value = -23349.5
$ws.Range("H100").Value = 2335.8333
$ws.Range("I100").Value = 2011.4286
$ws.Range("J100").Value = 2790
$ws.Range("K100").Value = 2011.4286
$ws.Range("L100").Value = 2790
$ws.Range("M100").Value = -1470.4286
$ws.Range("N100").Value = -3872
$ws.Range("H113").Value = 4126
$ws.Range("I113").Value = 3752
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 3752
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -1582
$ws.Range("N113").Value = -8840
$ws.Range("H132").Value = 1477.0927
$ws.Range("I132").Value = 1082.5
$ws.Range("J132").Value = 2858.1667
$ws.Range("K132").Value = 3247.5
$ws.Range("L132").Value = 8574.500100000001
$ws.Range("M132").Value = -717.5
$ws.Range("N132").Value = -13634.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 571.3570999999999
$ws.Range("I113").Value = 557.61536
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 1672.84608
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 497.15392
$ws.Range("N113").Value = -6590
$ws.Range("H132").Value = 1883.579
$ws.Range("I132").Value = 1780.8
$ws.Range("K132").Value = 5342.4
$ws.Range("M132").Value = -2812.4
